# Applies the "horizon to model area/shadow" data update to Sheet1.
# Only the numeric values in column D change; layout/formatting/formulas
# are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = 4058.788976559438
    "D3"  = 4058.788976559434
    "D4"  = 4058.788976559434
    "D6"  = 11596.53993302689
    "D7"  = 11596.53993302689
    "D9"  = 1856.573077663361
    "D10" = 1856.573077663359
    "D11" = 27
    "D12" = 1532.573077663359
    "D13" = 27
    "D17" = 21600
    "D19" = 15785.03079000041
    "D20" = 15785.03079000041
    "D21" = 324
    "D24" = 186936.8792052546
    "D25" = 186936.8792052545
    "D28" = 3738.737584105171
    "D29" = 3738.737584105171
    "D30" = 186936.8792052545
    "D35" = 20217.59999999949
    "D36" = 20217.59999999949
    "D38" = -352.4795252122792
    "D39" = -352.4795252122789
    "D41" = 352.4795252122789
    "D42" = 7049.590504245611
    "D43" = 7049.590504245611
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
